$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 322, pushing existing rows 322:405 down to 323:406
$ws.Rows.Item(322).Insert()

# Populate the newly inserted row 322 with the new record's data
$ws.Range("A322").Value = 3
$ws.Range("B322").Value = "Femacal de La Calera"
$ws.Range("C322").Value = "Coquimbo"
$ws.Range("D322").Value = "2022-03-22"
$ws.Range("E322").Value = 5
$ws.Range("F322").Value = 100112003
$ws.Range("G322").Value = "Ajo"
$ws.Range("H322").Value = "Chino"
$ws.Range("I322").Value = "Primera"
$ws.Range("J322").Value = 73
$ws.Range("K322").Value = 16500
$ws.Range("L322").Value = 17000
$ws.Range("M322").Value = 16740
$ws.Range("N322").Value = '$/caja 10 kilos'
$ws.Range("O322").Value = "China"
$ws.Range("P322").Value = 1674
$ws.Range("Q322").Value = 10
$ws.Range("R322").Value = "Hortaliza"
